$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = "9.01."
$ws.Range("B3").Value = "16.68."
$ws.Range("B4").Value = "16.01."
$ws.Range("B5").Value = "15.93."
$ws.Range("B6").Value = "15.37."
$ws.Range("B7").Value = "23.9."
